$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.034158668894515
$ws.Cells.Item(2, 4).Value = 1.03253858946127
$ws.Cells.Item(2, 5).Value = 1.042154106829229
$ws.Cells.Item(2, 6).Value = 1.050189332583513
$ws.Cells.Item(2, 9).Value = 1.034942598479826
$ws.Cells.Item(2, 10).Value = 1.039279444095624
$ws.Cells.Item(2, 11).Value = 1.035343530367846
$ws.Cells.Item(2, 12).Value = 1.044931596789768
$ws.Cells.Item(2, 13).Value = 1.05294430124736
$ws.Cells.Item(2, 14).Value = 1.040755340113258

$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.035106467916112
$ws.Cells.Item(3, 4).Value = 1.032995770851811
$ws.Cells.Item(3, 5).Value = 1.04303187717321
$ws.Cells.Item(3, 6).Value = 1.05123149412268
$ws.Cells.Item(3, 9).Value = 1.035088116697097
$ws.Cells.Item(3, 10).Value = 1.039870266532309
$ws.Cells.Item(3, 11).Value = 1.035610445813985
$ws.Cells.Item(3, 12).Value = 1.045619965182596
$ws.Cells.Item(3, 13).Value = 1.053798266992154
$ws.Cells.Item(3, 14).Value = 1.041347001585571

$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.035720204558769
$ws.Cells.Item(4, 4).Value = 1.033291860159364
$ws.Cells.Item(4, 5).Value = 1.04360067631121
$ws.Cells.Item(4, 6).Value = 1.051907007409007
$ws.Cells.Item(4, 9).Value = 1.035181225273844
$ws.Cells.Item(4, 10).Value = 1.040252402535821
$ws.Cells.Item(4, 11).Value = 1.035782692812457
$ws.Cells.Item(4, 12).Value = 1.046065570427277
$ws.Cells.Item(4, 13).Value = 1.054351399493403
$ws.Cells.Item(4, 14).Value = 1.041729680266041

$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.035978325654055
$ws.Cells.Item(5, 4).Value = 1.033416397214119
$ws.Cells.Item(5, 5).Value = 1.04383999530924
$ws.Cells.Item(5, 6).Value = 1.052191271002225
$ws.Cells.Item(5, 9).Value = 1.035220115825208
$ws.Cells.Item(5, 10).Value = 1.040413011988093
$ws.Cells.Item(5, 11).Value = 1.035854993048905
$ws.Cells.Item(5, 12).Value = 1.046252946175355
$ws.Cells.Item(5, 13).Value = 1.05458406933264
$ws.Cells.Item(5, 14).Value = 1.041890517802157

$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.036021671524717
$ws.Cells.Item(6, 4).Value = 1.033437311067446
$ws.Cells.Item(6, 5).Value = 1.043880189471555
$ws.Cells.Item(6, 6).Value = 1.052239016368788
$ws.Cells.Item(6, 9).Value = 1.035226630908482
$ws.Cells.Item(6, 10).Value = 1.040439976620164
$ws.Cells.Item(6, 11).Value = 1.035867125955003
$ws.Cells.Item(6, 12).Value = 1.046284409899595
$ws.Cells.Item(6, 13).Value = 1.054623143388157
$ws.Cells.Item(6, 14).Value = 1.041917520727098

$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.035723653169224
$ws.Cells.Item(7, 4).Value = 1.033293523991622
$ws.Cells.Item(7, 5).Value = 1.04360387333535
$ws.Cells.Item(7, 6).Value = 1.051910804662849
$ws.Cells.Item(7, 9).Value = 1.03518174592336
$ws.Cells.Item(7, 10).Value = 1.040254548767506
$ws.Cells.Item(7, 11).Value = 1.035783659334735
$ws.Cells.Item(7, 12).Value = 1.046068073982234
$ws.Cells.Item(7, 13).Value = 1.054354507918114
$ws.Cells.Item(7, 14).Value = 1.04173182954562

$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.034478888868874
$ws.Cells.Item(8, 4).Value = 1.032693040703355
$ws.Cells.Item(8, 5).Value = 1.042450582293843
$ws.Cells.Item(8, 6).Value = 1.050541294044928
$ws.Cells.Item(8, 9).Value = 1.034991994402995
$ws.Cells.Item(8, 10).Value = 1.039479149071199
$ws.Cells.Item(8, 11).Value = 1.03543383137576
$ws.Cells.Item(8, 12).Value = 1.045164195259572
$ws.Cells.Item(8, 13).Value = 1.053232786491335
$ws.Cells.Item(8, 14).Value = 1.040955328692805

$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.032288919481409
$ws.Cells.Item(9, 4).Value = 1.031636997044998
$ws.Cells.Item(9, 5).Value = 1.040424692327622
$ws.Cells.Item(9, 6).Value = 1.04813701381078
$ws.Cells.Item(9, 9).Value = 1.034649600728883
$ws.Cells.Item(9, 10).Value = 1.038111565863527
$ws.Cells.Item(9, 11).Value = 1.034813873497666
$ws.Cells.Item(9, 12).Value = 1.043572903428602
$ws.Cells.Item(9, 13).Value = 1.051260500853252
$ws.Cells.Item(9, 14).Value = 1.039585803360114

$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.030831319767605
$ws.Cells.Item(10, 4).Value = 1.030934470973479
$ws.Cells.Item(10, 5).Value = 1.039078442860777
$ws.Cells.Item(10, 6).Value = 1.046540258221745
$ws.Cells.Item(10, 9).Value = 1.034415975477514
$ws.Cells.Item(10, 10).Value = 1.037199068099276
$ws.Cells.Item(10, 11).Value = 1.03439826625418
$ws.Cells.Item(10, 12).Value = 1.042513078458737
$ws.Cells.Item(10, 13).Value = 1.049948615615427
$ws.Cells.Item(10, 14).Value = 1.038672009744373

$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.030200737700882
$ws.Cells.Item(11, 4).Value = 1.030630647397151
$ws.Cells.Item(11, 5).Value = 1.038496546830176
$ws.Cells.Item(11, 6).Value = 1.045850305168457
$ws.Cells.Item(11, 9).Value = 1.034313548343408
$ws.Cells.Item(11, 10).Value = 1.036803773087861
$ws.Cells.Item(11, 11).Value = 1.034217770839884
$ws.Cells.Item(11, 12).Value = 1.042054419878394
$ws.Cells.Item(11, 13).Value = 1.049381272048788
$ws.Cells.Item(11, 14).Value = 1.038276153368701

$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.029966597580366
$ws.Cells.Item(12, 4).Value = 1.030517851855199
$ws.Cells.Item(12, 5).Value = 1.038280561957105
$ws.Cells.Item(12, 6).Value = 1.045594245168922
$ws.Cells.Item(12, 9).Value = 1.034275312720648
$ws.Cells.Item(12, 10).Value = 1.036656917074018
$ws.Cells.Item(12, 11).Value = 1.034150647520718
$ws.Cells.Item(12, 12).Value = 1.041884092510012
$ws.Cells.Item(12, 13).Value = 1.049170643221751
$ws.Cells.Item(12, 14).Value = 1.038129088802474

$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.030016817520845
$ws.Cells.Item(13, 4).Value = 1.030542044226468
$ws.Cells.Item(13, 5).Value = 1.038326884312577
$ws.Cells.Item(13, 6).Value = 1.045649160968205
$ws.Cells.Item(13, 9).Value = 1.034283522963494
$ws.Cells.Item(13, 10).Value = 1.036688419343984
$ws.Cells.Item(13, 11).Value = 1.034165049271519
$ws.Cells.Item(13, 12).Value = 1.041920626536257
$ws.Cells.Item(13, 13).Value = 1.049215818925032
$ws.Cells.Item(13, 14).Value = 1.038160635809276

$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.030181381843885
$ws.Cells.Item(14, 4).Value = 1.030621322487953
$ws.Cells.Item(14, 5).Value = 1.0384786902393
$ws.Cells.Item(14, 6).Value = 1.045829134689646
$ws.Cells.Item(14, 9).Value = 1.034310391639965
$ws.Cells.Item(14, 10).Value = 1.036791634446474
$ws.Cells.Item(14, 11).Value = 1.034212224016713
$ws.Cells.Item(14, 12).Value = 1.042040339767267
$ws.Cells.Item(14, 13).Value = 1.049363859198908
$ws.Cells.Item(14, 14).Value = 1.038263997489052

$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.030282786791452
$ws.Cells.Item(15, 4).Value = 1.030670176187223
$ws.Cells.Item(15, 5).Value = 1.038572243748593
$ws.Cells.Item(15, 6).Value = 1.045940051527831
$ws.Cells.Item(15, 9).Value = 1.034326921207879
$ws.Cells.Item(15, 10).Value = 1.036855225277207
$ws.Cells.Item(15, 11).Value = 1.03424127945948
$ws.Cells.Item(15, 12).Value = 1.042114104214791
$ws.Cells.Item(15, 13).Value = 1.04945508601535
$ws.Cells.Item(15, 14).Value = 1.038327678626057

$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.030873181459757
$ws.Cells.Item(16, 4).Value = 1.030954642790398
$ws.Cells.Item(16, 5).Value = 1.039117083369245
$ws.Cells.Item(16, 6).Value = 1.046586078827983
$ws.Cells.Item(16, 9).Value = 1.034422746606353
$ws.Cells.Item(16, 10).Value = 1.037225298857988
$ws.Cells.Item(16, 11).Value = 1.034410233964006
$ws.Cells.Item(16, 12).Value = 1.042543523517112
$ws.Cells.Item(16, 13).Value = 1.049986283409224
$ws.Cells.Item(16, 14).Value = 1.038698277753771

$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.031243673106915
$ws.Cells.Item(17, 4).Value = 1.031133182767089
$ws.Cells.Item(17, 5).Value = 1.039459125810577
$ws.Cells.Item(17, 6).Value = 1.04699170433491
$ws.Cells.Item(17, 9).Value = 1.034482516833332
$ws.Cells.Item(17, 10).Value = 1.037457389302134
$ws.Cells.Item(17, 11).Value = 1.034516072240897
$ws.Cells.Item(17, 12).Value = 1.042812955167164
$ws.Cells.Item(17, 13).Value = 1.050319680634998
$ws.Cells.Item(17, 14).Value = 1.038930697792969

$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.031459829390582
$ws.Cells.Item(18, 4).Value = 1.03123735816565
$ws.Cells.Item(18, 5).Value = 1.039658733539082
$ws.Cells.Item(18, 6).Value = 1.047228438920507
$ws.Cells.Item(18, 9).Value = 1.034517257621139
$ws.Cells.Item(18, 10).Value = 1.037592746554652
$ws.Cells.Item(18, 11).Value = 1.034577754252565
$ws.Cells.Item(18, 12).Value = 1.042970134404679
$ws.Cells.Item(18, 13).Value = 1.050514214171546
$ws.Cells.Item(18, 14).Value = 1.039066247268311

$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.031533542364752
$ws.Cells.Item(19, 4).Value = 1.031272885346064
$ws.Cells.Item(19, 5).Value = 1.039726811570207
$ws.Cells.Item(19, 6).Value = 1.047309183026558
$ws.Cells.Item(19, 9).Value = 1.034529082592196
$ws.Cells.Item(19, 10).Value = 1.037638896941556
$ws.Cells.Item(19, 11).Value = 1.03459877742927
$ws.Cells.Item(19, 12).Value = 1.043023732559727
$ws.Cells.Item(19, 13).Value = 1.050580556682887
$ws.Cells.Item(19, 14).Value = 1.039112463194058

$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.031203917188356
$ws.Cells.Item(20, 4).Value = 1.031114023375769
$ws.Cells.Item(20, 5).Value = 1.039422417525446
$ws.Cells.Item(20, 6).Value = 1.046948170059103
$ws.Cells.Item(20, 9).Value = 1.034476116686684
$ws.Cells.Item(20, 10).Value = 1.037432489964736
$ws.Cells.Item(20, 11).Value = 1.034504722131743
$ws.Cells.Item(20, 12).Value = 1.042784045179534
$ws.Cells.Item(20, 13).Value = 1.050283903167911
$ws.Cells.Item(20, 14).Value = 1.038905763095656

$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.030132919368461
$ws.Cells.Item(21, 4).Value = 1.03059797539913
$ws.Cells.Item(21, 5).Value = 1.038433982828269
$ws.Cells.Item(21, 6).Value = 1.045776130852707
$ws.Cells.Item(21, 9).Value = 1.034302484714523
$ws.Cells.Item(21, 10).Value = 1.036761240878283
$ws.Cells.Item(21, 11).Value = 1.034198334410293
$ws.Cells.Item(21, 12).Value = 1.042005086124846
$ws.Cells.Item(21, 13).Value = 1.049320262043608
$ws.Cells.Item(21, 14).Value = 1.038233560758506

$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.029460039269478
$ws.Cells.Item(22, 4).Value = 1.030273852341809
$ws.Cells.Item(22, 5).Value = 1.037813424906361
$ws.Cells.Item(22, 6).Value = 1.045040492516793
$ws.Cells.Item(22, 9).Value = 1.034192218568888
$ws.Cells.Item(22, 10).Value = 1.036339050803395
$ws.Cells.Item(22, 11).Value = 1.034005238085112
$ws.Cells.Item(22, 12).Value = 1.04151554887063
$ws.Cells.Item(22, 13).Value = 1.048715006995119
$ws.Cells.Item(22, 14).Value = 1.037810771125286

$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.029816698209404
$ws.Cells.Item(23, 4).Value = 1.030445643602763
$ws.Cells.Item(23, 5).Value = 1.038142307743232
$ws.Cells.Item(23, 6).Value = 1.04543034766055
$ws.Cells.Item(23, 9).Value = 1.034250776520107
$ws.Cells.Item(23, 10).Value = 1.036562875613521
$ws.Cells.Item(23, 11).Value = 1.034107645182524
$ws.Cells.Item(23, 12).Value = 1.04177504018333
$ws.Cells.Item(23, 13).Value = 1.049035804644905
$ws.Cells.Item(23, 14).Value = 1.038034913792316

$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.03122188100746
$ws.Cells.Item(24, 4).Value = 1.031122680568836
$ws.Cells.Item(24, 5).Value = 1.03943900411003
$ws.Cells.Item(24, 6).Value = 1.046967840891643
$ws.Cells.Item(24, 9).Value = 1.03447900901502
$ws.Cells.Item(24, 10).Value = 1.037443740956658
$ws.Cells.Item(24, 11).Value = 1.034509850917275
$ws.Cells.Item(24, 12).Value = 1.042797108283467
$ws.Cells.Item(24, 13).Value = 1.050300069253442
$ws.Cells.Item(24, 14).Value = 1.038917030065277

$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.032854663318857
$ws.Cells.Item(25, 4).Value = 1.031909752126236
$ws.Cells.Item(25, 5).Value = 1.040947672965847
$ws.Cells.Item(25, 6).Value = 1.048757507883884
$ws.Cells.Item(25, 9).Value = 1.034739065443679
$ws.Cells.Item(25, 10).Value = 1.03846525937403
$ws.Cells.Item(25, 11).Value = 1.034974557753585
$ws.Cells.Item(25, 12).Value = 1.043984111971827
$ws.Cells.Item(25, 13).Value = 1.051769864565508
$ws.Cells.Item(25, 14).Value = 1.039939999155972
